$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.561171122214548
$ws.Range("B1").Value = -0.0002517101720076429
$ws.Range("C1").Value = 1.400073368420803
$ws.Range("D1").Value = 0.1709741328928501
$ws.Range("E1").Value = 1.570796365408233
$ws.Range("F1").Value = 0.9903747278261085
$ws.Range("A2").Value = 2.555848965682979
$ws.Range("B2").Value = 0.0001803432818697641
$ws.Range("C2").Value = 1.398460229945526
$ws.Range("D2").Value = 0.1721552197852622
$ws.Range("E2").Value = 1.570796365220714
$ws.Range("F2").Value = 0.9850525714503057
$ws.Range("A3").Value = 2.522657562997278
$ws.Range("B3").Value = 0.002874826257770738
$ws.Range("C3").Value = 1.388399960503409
$ws.Range("D3").Value = 0.1795210179199698
$ws.Range("E3").Value = 1.570796364051259
$ws.Range("F3").Value = 0.9518611697360345
$ws.Range("A4").Value = 2.444161474445256
$ws.Range("B4").Value = 0.00924715032858687
$ws.Range("C4").Value = 1.364607901804851
$ws.Range("D4").Value = 0.1969407801431379
$ws.Range("E4").Value = 1.570796361285554
$ws.Range("F4").Value = 0.873365083481396
$ws.Range("A5").Value = 2.312809919049192
$ws.Range("B5").Value = 0.01991028901013736
$ws.Range("C5").Value = 1.324795422307113
$ws.Range("D5").Value = 0.2260901671361768
$ws.Range("E5").Value = 1.570796356657557
$ws.Range("F5").Value = 0.7420135319296621
$ws.Range("A6").Value = 2.129289331443423
$ws.Range("B6").Value = 0.03480851743734772
$ws.Range("C6").Value = 1.269170579423171
$ws.Range("D6").Value = 0.2668168461098673
$ws.Range("E6").Value = 1.570796350191454
$ws.Range("F6").Value = 0.558492949695081
$ws.Range("A7").Value = 1.900875918751935
$ws.Range("B7").Value = 0.05335115204042857
$ws.Range("C7").Value = 1.199938781730575
$ws.Range("D7").Value = 0.3175060894984869
$ws.Range("E7").Value = 1.570796342143612
$ws.Range("F7").Value = 0.3300795436886814
$ws.Range("A8").Value = 1.639788217465948
$ws.Range("B8").Value = 0.07454629022105426
$ws.Range("C8").Value = 1.120803451180303
$ws.Range("D8").Value = 0.3754463736539357
$ws.Range("E8").Value = 1.570796332944535
$ws.Range("F8").Value = 0.06899185004407737
$ws.Range("A9").Value = 1.361539650321506
$ws.Range("B9").Value = 0.09713455002854179
$ws.Range("C9").Value = 1.036466685305622
$ws.Range("D9").Value = 0.437194977539862
$ws.Range("E9").Value = 1.570796323140818
$ws.Range("F9").Value = -0.2092567089567264
$ws.Range("A10").Value = 1.083291083177064
$ws.Range("B10").Value = 0.1197228098360293
$ws.Range("C10").Value = 0.9521299194309416
$ws.Range("D10").Value = 0.4989435814257883
$ws.Range("E10").Value = 1.570796313337101
$ws.Range("F10").Value = -0.4875052679575304
$ws.Range("A11").Value = 0.8222033818910772
$ws.Range("B11").Value = 0.140917948016655
$ws.Range("C11").Value = 0.8729945888806701
$ws.Range("D11").Value = 0.5568838655812371
$ws.Range("E11").Value = 1.570796304138024
$ws.Range("F11").Value = -0.748592961602134
$ws.Range("A12").Value = 0.5937899691995894
$ws.Range("B12").Value = 0.1594605826197359
$ws.Range("C12").Value = 0.8037627911880737
$ws.Range("D12").Value = 0.6075731089698566
$ws.Range("E12").Value = 1.570796296090182
$ws.Range("F12").Value = -0.977006367608534
$ws.Range("A13").Value = 0.410269381593821
$ws.Range("B13").Value = 0.1743588110469463
$ws.Range("C13").Value = 0.7481379483041316
$ws.Range("D13").Value = 0.6482997879435473
$ws.Range("E13").Value = 1.570796289624078
$ws.Range("F13").Value = -1.160526949843115
$ws.Range("A14").Value = 0.2789178261977563
$ws.Range("B14").Value = 0.1850219497284968
$ws.Range("C14").Value = 0.7083254688063935
$ws.Range("D14").Value = 0.6774491749365862
$ws.Range("E14").Value = 1.570796284996081
$ws.Range("F14").Value = -1.291878501394849
$ws.Range("A15").Value = 0.2004217376457329
$ws.Range("B15").Value = 0.1913942737993129
$ws.Range("C15").Value = 0.6845334101078362
$ws.Range("D15").Value = 0.6948689371597545
$ws.Range("E15").Value = 1.570796282230376
$ws.Range("F15").Value = -1.370374587649487
$ws.Range("A16").Value = 0.1672303349600339
$ws.Range("B16").Value = 0.194088756775214
$ws.Range("C16").Value = 0.6744731406657185
$ws.Range("D16").Value = 0.7022347352944619
$ws.Range("E16").Value = 1.570796281060922
$ws.Range("F16").Value = -1.403565989363758
$ws.Range("A17").Value = 0.1619081784284635
$ws.Range("B17").Value = 0.1945208102290914
$ws.Range("C17").Value = 0.6728600021904421
$ws.Range("D17").Value = 0.7034158221868742
$ws.Range("E17").Value = 1.570796280873402
$ws.Range("F17").Value = -1.40888814573956
